$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# Row 2: the Approved/Rejected result moves to "Rejected", with a reason
# ("tet") added in the ReasonToReject column.
$ws.Range("I2").Value = "Rejected"
$ws.Range("J2").Value = "tet"

# Row 19: was "Rejected" with reason "test" -- now "Approved" with no
# reason (cleared).
$ws.Range("I19").Value = "Approved"
$ws.Range("J19").ClearContents()

# Update the active selection to J19.
$ws.Range("J19").Select()
